# Update column G ("K" - strikeouts) values on Sheet1 for rows 2-33.
# These values were regenerated from source data (switching from "Strike#" to "K"),
# so we simply overwrite the previous numbers with the new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 5
    3  = 1
    4  = 4
    5  = 2
    6  = 3
    7  = 5
    8  = 3
    9  = 3
    10 = 4
    11 = 3
    12 = 3
    13 = 2
    14 = 4
    15 = 3
    16 = 3
    17 = 2
    18 = 1
    19 = 1
    20 = 0
    21 = 4
    22 = 1
    23 = 2
    24 = 0
    25 = 5
    26 = 0
    27 = 3
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 3
    33 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
